$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 (copy formatting from neighboring header cell G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the per-row "Save" indicator values (H2:H41)
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(6, 8).Value = 1
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(11, 8).Value = 1
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(13, 8).Value = 1
$ws.Cells.Item(14, 8).Value = 1
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(21, 8).Value = 1
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(27, 8).Value = 1
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(40, 8).Value = 1
$ws.Cells.Item(41, 8).Value = 0
